$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 22 (shifts existing rows 22-23 down to 23-24)
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row with the Credit sheet name entry
$ws.Range("A22").Value = "CreditSheetName"
$ws.Range("B22").Value = "Sheet1"

# Keep the active selection in sync with where Excel would leave the cursor
$ws.Range("B26").Select()
